$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register a Bold font in the workbook's font table (styles.xml <fonts>) ---
# Touching Font.Bold on the C10 cell (before we overwrite its value) forces the
# engine to add a "<b/>..." font entry, matching the font the real edit added
# when the author bolded the "git push -f origin main" run.
$regCell = $ws.Cells.Item(10, 3)
$regCell.Font.Bold = $true
$regCell.Font.Bold = $false

# --- Update the "Notes" cell (C10) for the "git push -u origin main" row ---
# Replace the old single-run note with a two-run rich text note: a plain
# paragraph followed by a bold command line.
$part1 = "There may be times where you aren't able to push local work because the work on GitHub has been changed also. The solution to this is to pull first, sort out the errors, then push once again.`nIf you really want to force through a change (could be risky if there's multiple people), do this command:`n"
$part2 = "git push -f origin main"

$cellC10 = $ws.Cells.Item(10, 3)
$cellC10.Value2 = $part1 + $part2

$boldChars = $cellC10.Characters($part1.Length + 1, $part2.Length)
$boldChars.Font.Bold = $true

# --- Row height updates ---
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 105

# --- Selection moved to A10 ---
$ws.Range("A10").Select()
